# Trade #208 closed at 2026-02-17 16:49:19 - unknown UNKNOWN +0.000%
#
# Updates the Summary, Strategy Status, All Trades, volatility_scorer and
# MarketMaking sheets to reflect:
#   - Trade #208 (volatility_scorer, CLOSED) being appended
#   - Trade #209 (MarketMaking, OPEN) being appended
#   - The roll-up stats on Summary / Strategy Status recomputed accordingly

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.17    # Current Capital
$summary.Range("B4").Value = -0.84      # Total P&L $
$summary.Range("B6").Value = 208        # Total Trades
$summary.Range("B8").Value = 100        # Losing Trades
$summary.Range("B9").Value = 41.35      # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - volatility_scorer row (row 12)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C12").Value = 99.93000000000001  # Capital
$status.Range("D12").Value = 11                 # Trades
$status.Range("E12").Value = -0.08              # P&L $
$status.Range("F12").Value = -0.07000000000000001  # P&L %
$status.Range("G12").Value = 36.36              # Win Rate %

# ---------------------------------------------------------------------------
# Helper: write one trade-log row (shared shape across the trade sheets)
# Columns: A # | B Date | C Time | D Strategy | E Side | F Entry | G Exit |
#          H Status | I P&L % | J P&L $ | K Capital | L ? | M ? | N ? |
#          O Note | P Exit Reason | Q ?
#
# NOTE: positional parameters only - this COM/PowerShell host does not wire
# up named (`-paramName value`) arguments correctly.
# ---------------------------------------------------------------------------
function Write-TradeRow($sheet, $row, $tradeNum, $date, $time, $strategy, $side, $entry, $exit, $status, $pnlPct, $pnlDollar, $capital, $l, $m, $n, $note, $exitReason, $q) {

    $sheet.Range("A$row").Value = $tradeNum

    # Force the date/time columns to be stored as literal text (matching the
    # existing rows) rather than letting Excel auto-coerce them into date /
    # time serial numbers.
    $sheet.Range("B$row").NumberFormat = "@"
    $sheet.Range("B$row").Value = $date
    $sheet.Range("B$row").Style = "Normal"

    $sheet.Range("C$row").Value = $time
    $sheet.Range("D$row").Value = $strategy
    $sheet.Range("E$row").Value = $side

    $sheet.Range("F$row").Value = $entry

    if ($null -eq $exit) {
        $sheet.Range("G$row").Value = ""
    } else {
        $sheet.Range("G$row").Value = $exit
    }

    $sheet.Range("H$row").Value = $status
    $sheet.Range("I$row").Value = $pnlPct
    $sheet.Range("J$row").Value = $pnlDollar
    $sheet.Range("K$row").Value = $capital
    $sheet.Range("L$row").Value = $l
    $sheet.Range("M$row").Value = $m
    $sheet.Range("N$row").Value = $n
    $sheet.Range("O$row").Value = $note

    if ($null -eq $exitReason) {
        $sheet.Range("P$row").Value = ""
    } else {
        $sheet.Range("P$row").Value = $exitReason
    }

    $sheet.Range("Q$row").Value = $q
}

# ---------------------------------------------------------------------------
# All Trades sheet - append trade #208 (row 209) and trade #209 (row 210)
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# row, tradeNum, date, time, strategy, side, entry, exit, status, pnlPct, pnlDollar, capital, l, m, n, note, exitReason, q
Write-TradeRow $allTrades 209 208 "2026-02-17" "16:49:12" "volatility_scorer" "NEUTRAL" 0.24 0.193089 "CLOSED" -19.5461 -0.05 99.93000000000001 0 0 0.85 "Low vol market (score: inf) - ideal for market making" "early_exit" 0.18

Write-TradeRow $allTrades 210 209 "2026-02-17" "16:49:13" "MarketMaking" "UP" 0.77 $null "OPEN" 0 0 99.24059897733886 0 0 0.6 "Normal spread capture: 19600 bps" $null 0

# ---------------------------------------------------------------------------
# volatility_scorer sheet - append trade #208 (row 12)
# ---------------------------------------------------------------------------
$volScorer = $wb.Worksheets.Item("volatility_scorer")

Write-TradeRow $volScorer 12 208 "2026-02-17" "16:49:12" "volatility_scorer" "NEUTRAL" 0.24 0.193089 "CLOSED" -19.5461 -0.05 99.93000000000001 0 0 0.85 "Low vol market (score: inf) - ideal for market making" "early_exit" 0.18

# ---------------------------------------------------------------------------
# MarketMaking sheet - append trade #209 (row 199)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

Write-TradeRow $marketMaking 199 209 "2026-02-17" "16:49:13" "MarketMaking" "UP" 0.77 $null "OPEN" 0 0 99.24059897733886 0 0 0.6 "Normal spread capture: 19600 bps" $null 0
